$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2023-06-29 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-06-30 Friday", 2)

# Replace every arithmetic-answer cell in the table (20 rows x 5 columns,
# row-major order) with its new value. Cell.Range.Text is used instead of
# Find/Replace so duplicate old values (e.g. two cells both reading
# "45+15=60") are disambiguated correctly and resolve to the distinct
# replacement defined for that specific cell position.
$newValues = @(
    "68-42=26",
    "68-25=43",
    "1+23=24",
    "70+0=70",
    "62-36=26",
    "27+48=75",
    "45-38=7",
    "50-47=3",
    "44+54=98",
    "9+15=24",
    "5+75=80",
    "49-3=46",
    "20+7=27",
    "70+17=87",
    "40-26=14",
    "83-51=32",
    "58+25=83",
    "89-18=71",
    "12-2=10",
    "51+45=96",
    "70-64=6",
    "13+85=98",
    "46-45=1",
    "41-21=20",
    "69-13=56",
    "44+6=50",
    "61-6=55",
    "56+3=59",
    "73+11=84",
    "79+17=96",
    "95-16=79",
    "6+1=7",
    "52-6=46",
    "55-29=26",
    "92-55=37",
    "45+52=97",
    "37+30=67",
    "69-17=52",
    "56-5=51",
    "89-15=74",
    "3+12=15",
    "30+49=79",
    "5+26=31",
    "0+65=65",
    "51+2=53",
    "75+18=93",
    "75-21=54",
    "38+8=46",
    "4+64=68",
    "26+27=53",
    "83+5=88",
    "81+15=96",
    "66-22=44",
    "20+75=95",
    "69+11=80",
    "13+9=22",
    "55-49=6",
    "39+30=69",
    "61+2=63",
    "51-21=30",
    "4+31=35",
    "6+33=39",
    "49+43=92",
    "20-13=7",
    "45+39=84",
    "82-18=64",
    "25+1=26",
    "81-73=8",
    "38+15=53",
    "1+54=55",
    "56-40=16",
    "71-1=70",
    "6+74=80",
    "18+75=93",
    "63+26=89",
    "98-1=97",
    "66-13=53",
    "85-6=79",
    "54+9=63",
    "49+16=65",
    "81-73=8",
    "70+9=79",
    "48+20=68",
    "25+23=48",
    "52+22=74",
    "60-9=51",
    "83+7=90",
    "53+5=58",
    "46+0=46",
    "20+54=74",
    "53-11=42",
    "92-64=28",
    "47+38=85",
    "47+36=83",
    "35+25=60",
    "41+56=97",
    "70-60=10",
    "36+48=84",
    "58-11=47",
    "65-23=42"
)

$t = $d.Tables.Item(1)
$i = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$i]
        $i = $i + 1
    }
}
